# Daily attendance processing - 2025-11-28 17:49:03
# Reorder the "Recorded By" (column G) comma-separated list so that any
# "System"/"system" entries are moved to the front, preserving the
# relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq $null) { continue }

    $parts = $val.Split(",")
    if ($parts.Length -le 1) { continue }

    $sysList = @()
    $otherList = @()

    foreach ($p in $parts) {
        $t = $p.Trim()
        if ($t.ToLower() -eq "system") {
            $sysList += $t
        } else {
            $otherList += $t
        }
    }

    if ($sysList.Length -eq 0) { continue }

    $combined = $sysList + $otherList
    $newVal = $combined -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
